# Fix Arduino bill of materials
#
# 1. Slide 3: "Should total no more than $60 dollars overall (pre-tax)"
#    -> "Should total no more than $40 dollars overall (pre-tax)"
#    (budget lowered now that the USB-to-TTL converter is no longer a
#    separate purchase)
#
# 2. Slide 38 ("Building your own"): drop the standalone "USB to TTL
#    converter ($20 on Amazon)" line item and the related "cheaper parts /
#    USB -> TTL" commentary line, since they're no longer part of the
#    bill of materials.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3: budget line $60 -> $40
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$shBudget = $s3.Shapes.Item(2)
$trBudget = $shBudget.TextFrame.TextRange

$budgetParaCount = $trBudget.Paragraphs().Count
for ($i = 1; $i -le $budgetParaCount; $i++) {
    $para = $trBudget.Paragraphs($i, 1)
    if ($para.Text -like "*`$60 dollars overall*") {
        $startInPara = $para.Text.IndexOf("`$60 ") + 1
        $sub = $para.Characters($startInPara, 4)
        $sub.Text = "`$40 "
        break
    }
}

# ---------------------------------------------------------------------
# Slide 38: Arduino parts list cleanup
# ---------------------------------------------------------------------
$s38 = $p.Slides.Item(38)
$shParts = $s38.Shapes.Item(2)
$trParts = $shParts.TextFrame.TextRange

# Remove the "USB to TTL converter ($20 on Amazon)" paragraph entirely.
$count = $trParts.Paragraphs().Count
for ($i = 1; $i -le $count; $i++) {
    $para = $trParts.Paragraphs($i, 1)
    if ($para.Text -like "USB to TTL converter*") {
        $para.Delete()
        break
    }
}

# Split "Battery with female JST connector (" into two runs:
# "Battery " + "with female JST connector ("
$count = $trParts.Paragraphs().Count
for ($i = 1; $i -le $count; $i++) {
    $para = $trParts.Paragraphs($i, 1)
    if ($para.Text -like "Battery with female JST connector*") {
        $sub = $para.Characters(1, 8)
        $sub.Text = "Battery "
        break
    }
}

# Remove the "You could conceivably get cheaper parts..." paragraph entirely.
$count = $trParts.Paragraphs().Count
for ($i = 1; $i -le $count; $i++) {
    $para = $trParts.Paragraphs($i, 1)
    if ($para.Text -like "You could conceivably*") {
        $para.Delete()
        break
    }
}

# Split "These are the ones I know work, using other parts YMMV" into two
# runs: "These " + "are the ones I know work, using other parts YMMV"
$count = $trParts.Paragraphs().Count
for ($i = 1; $i -le $count; $i++) {
    $para = $trParts.Paragraphs($i, 1)
    if ($para.Text -like "These are the ones*") {
        $sub = $para.Characters(1, 6)
        $sub.Text = "These "
        break
    }
}
